# Updated cryptos list — refresh Price (col D) and Volume(1h) (col E) figures,
# and fix the ImmutableX / Bittensor row ordering (row 34 <-> row 35, each also
# getting a refreshed Price/Volume).
#
# Price cells are plain text (dotted thousands separators, e.g. "63.009.94"),
# and Volume cells are padded percentage strings (e.g. "  -0.95%  "). Both
# must round-trip as literal text: pre-setting NumberFormat to "@" (Text)
# before the assignment stops Excel's COM layer from re-interpreting
# numeric-looking strings (e.g. "582.39", "18.30", "0.0538") as real numbers
# and mangling trailing zeros / switching to scientific notation. Resetting
# the Style back to "Normal" afterwards keeps the cell's style index
# identical to the original (no lingering explicit number format).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2"  = "63.009.94"
    "E2"  = "  -0.95%  "
    "D3"  = "2.548.04"
    "E3"  = "  -0.12%  "
    "E4"  = "  -0.01%  "
    "D5"  = "582.39"
    "E5"  = "  +2.03%  "
    "D6"  = "146.75"
    "E6"  = "  -2.96%  "
    "E7"  = "  +0.01%  "
    "E8"  = "  -0.55%  "
    "E9"  = "  -0.69%  "
    "E10" = "  -3.65%  "
    "E11" = "  -0.17%  "
    "E12" = "  -1.53%  "
    "D13" = "27.48"
    "E13" = "  -3.68%  "
    "E14" = "  -0.02%  "
    "D15" = "62.927.48"
    "E15" = "  -0.91%  "
    "E16" = "  -1.22%  "
    "E17" = "  +0.12%  "
    "D18" = "11.33"
    "E18" = "  -3.29%  "
    "D19" = "338.45"
    "E19" = "  -0.84%  "
    "D20" = "4.32"
    "D21" = "6.75"
    "E22" = "  -0.13%  "
    "D23" = "65.62"
    "E23" = "  -0.81%  "
    "D24" = "2.678.63"
    "E24" = "  +0.32%  "
    "D26" = "1.60"
    "E26" = "  -0.55%  "
    "E27" = "  -3.44%  "
    "E28" = "  +0.06%  "
    "E29" = "  -3.32%  "
    "D30" = "7.66"
    "E30" = "  +5.75%  "
    "E31" = "  +4.25%  "
    "D32" = "0.0₃0814"
    "E32" = "  -2.43%  "
    "D33" = "178.07"
    "E33" = "  -0.06%  "
    # row 34 was ImmutableX -> now Bittensor (moved up from row 35)
    "B34" = "Bittensor"
    "C34" = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
    "D34" = "417.51"
    "E34" = "  -1.69%  "
    # row 35 was Bittensor -> now ImmutableX (moved down from row 34)
    "B35" = "ImmutableX"
    "C35" = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
    "D35" = "1.54"
    "E35" = "  -2.70%  "
    "E36" = "  -1.61%  "
    "D37" = "19.10"
    "E37" = "  -0.61%  "
    "E39" = "  -2.66%  "
    "E40" = "  -2.91%  "
    "E41" = "  +0.04%  "
    "D42" = "39.77"
    "E42" = "  -0.02%  "
    "D43" = "150.72"
    "E43" = "  -2.19%  "
    "E44" = "  -1.33%  "
    "D45" = "20.74"
    "E45" = "  -1.96%  "
    "D46" = "0.0538"
    "E46" = "  +1.08%  "
    "E47" = "  -1.57%  "
    "E48" = "  -0.14%  "
    "D49" = "0.0237"
    "E49" = "  -0.91%  "
    "D50" = "18.30"
    "E50" = "  -2.23%  "
    "D51" = "1.70"
    "E51" = "  -6.98%  "
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}
